$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "datos actualizados" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 7 de Octubre de 2020 a las 21:21"

# Row 4
$ws.Range("B4").Value = 7753872
$ws.Range("C4").Value = 30720
$ws.Range("D4").Value = 4963229
$ws.Range("E4").Value = 2574283
$ws.Range("G4").Value = 538
$ws.Range("H4").Value = 216360

# Row 5
$ws.Range("B5").Value = 6832646
$ws.Range("C5").Value = 78467
$ws.Range("D5").Value = 5821423
$ws.Range("E5").Value = 905672
$ws.Range("G5").Value = 960
$ws.Range("H5").Value = 105551

# Row 14
$ws.Range("D14").Value = 99793
$ws.Range("E14").Value = 521271

# Row 26
$ws.Range("B26").Value = 310955
$ws.Range("C26").Value = 3836
$ws.Range("E26").Value = 33603
$ws.Range("G26").Value = 17
$ws.Range("H26").Value = 9652

# Row 33
$ws.Range("A33").Value = "Marruecos"
$ws.Range("B33").Value = 140024
$ws.Range("C33").Value = 2776
$ws.Range("D33").Value = 118142
$ws.Range("E33").Value = 19443
$ws.Range("G33").Value = 29
$ws.Range("H33").Value = 2439

# Row 34
$ws.Range("A34").Value = "Bolivia"
$ws.Range("B34").Value = 137468
$ws.Range("C34").Value = 361
$ws.Range("D34").Value = 98542
$ws.Range("E34").Value = 30770
$ws.Range("G34").Value = 27
$ws.Range("H34").Value = 8156

# Row 72
$ws.Range("B72").Value = 39907
$ws.Range("C72").Value = 321
$ws.Range("E72").Value = 11828
$ws.Range("G72").Value = 5
$ws.Range("H72").Value = 748

# Row 93
$ws.Range("A93").Value = "Zambia"
$ws.Range("B93").Value = 15224
$ws.Range("C93").Value = 54
$ws.Range("D93").Value = 14342
$ws.Range("E93").Value = 547
$ws.Range("G93").Value = 0
$ws.Range("H93").Value = 335

# Row 94
$ws.Range("A94").Value = "Senegal"
$ws.Range("B94").Value = 15174
$ws.Range("C94").Value = 33
$ws.Range("D94").Value = 12998
$ws.Range("E94").Value = 1863
$ws.Range("G94").Value = 1
$ws.Range("H94").Value = 313

# Row 106
$ws.Range("B106").Value = 10103
$ws.Range("C106").Value = 33
$ws.Range("D106").Value = 9773
$ws.Range("E106").Value = 261

# Row 109
$ws.Range("B109").Value = 9494
$ws.Range("C109").Value = 96
$ws.Range("D109").Value = 6812
$ws.Range("E109").Value = 2614
$ws.Range("G109").Value = 1
$ws.Range("H109").Value = 68

# Row 138
$ws.Range("B138").Value = 4133
$ws.Range("C138").Value = 25
$ws.Range("D138").Value = 3678
$ws.Range("E138").Value = 424

# Row 145
$ws.Range("B145").Value = 3210
$ws.Range("C145").Value = 15
$ws.Range("D145").Value = 2502
$ws.Range("E145").Value = 577

# Row 189
$ws.Range("B189").Value = 227
$ws.Range("C189").Value = 3
$ws.Range("D189").Value = 202
$ws.Range("E189").Value = 23

# Row 207
$ws.Range("A207").Value = "Nueva Caledonia"

# Row 208
$ws.Range("A208").Value = "Santa Lucia"

# Row 215
$ws.Range("A215").Value = "Montserrat"
$ws.Range("D215").Value = 12
$ws.Range("H215").Value = 1

# Row 216
$ws.Range("A216").Value = "Islas Malvinas"
$ws.Range("D216").Value = 13
$ws.Range("H216").Value = 0
